$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at row 11 (pushes old rows 11-16 down to 12-17),
# preserving the blank separator row pattern between asset groups.
$ws.Rows("11").Insert()

# Fill in the new "Button trophy" asset row (row 11).
$ws.Range("A11").Value = "Button trophy"
$ws.Range("C11").Value = "CC BY 3.0  "
$ws.Range("D11").Value = "CC BY 3.0  "
$ws.Range("E11").Value = "http://creativecommons.org/licenses/by/3.0/"
$ws.Range("F11").Value = "https://www.makeschool.com/academy/art/object/game-icons/trophy"
$ws.Range("C11").ClearFormats()

# Update the License Name column for the existing CC-BY rows (anthem, cog-lock,
# Button Shopping Cart) from the old "Attribution, Commercially, Modify" label
# to "CC BY 3.0", and drop their special cell style so they render like the
# rest of the CC BY 3.0 rows.
$ws.Range("C8").Value = "CC BY 3.0  "
$ws.Range("C8").ClearFormats()

$ws.Range("C9").Value = "CC BY 3.0  "
$ws.Range("C9").ClearFormats()

$ws.Range("C10").Value = "CC BY 3.0  "
$ws.Range("C10").ClearFormats()

# Update the selection to match the edited region.
$ws.Range("C8:C11").Select()
